$d = $word.ActiveDocument

# Locate the paragraph that ends the "Al zijn personeel verdient ..." sentence
$r = $d.Content
$found = $r.Find.Execute("Al zijn personeel verdient  11.5 euro per uur. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetPara = $r.Paragraphs(1)

# Insert a new paragraph right after it and fill it with "Console: "
$targetPara.Range.InsertParagraphAfter()
$consolePara = $targetPara.Next()
$consolePara.Range.Text = "Console: "

# Insert another new paragraph after that one with the Johan quote
$consolePara.Range.InsertParagraphAfter()
$quotePara = $consolePara.Next()
$quotePara.Range.Text = "`"Johan heeft op donderdag (en vrijdag) 6 uur gewerkt en heeft hiervoor 69 euro verdiend.`""
